# Regenerate the "K" (strikeouts) column (G) values for the save_data sheet.
# This mirrors the upstream change: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals" - the K column values were
# recomputed and rewritten for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 3
    12 = 0
    13 = 2
    14 = 2
    15 = 2
    16 = 0
    17 = 3
    18 = 1
    19 = 0
    20 = 2
    21 = 0
    22 = 1
    23 = 1
    24 = 1
    26 = 1
    27 = 2
    28 = 0
    30 = 1
    31 = 1
    33 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
